# Updated cryptos list on Tue Oct 17 17:41:16 UTC 2023 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) columns with latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.483.27'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = '1.572.13'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.99'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '46.31'
$ws.Range('E8').Value = '  +6.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '24.11'
$ws.Range('E9').Value = '  +2.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.247'
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0590'
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0882'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('D13').Value = '1.797.91'
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').Value = '1.570.29'
$ws.Range('E14').Value = '  -0.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.521'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.69'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('D17').Value = '28.503.22'
$ws.Range('E17').Value = '  +1.32%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.13'
$ws.Range('E18').Value = '  -2.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.36'
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.36'
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('D21').Value = '0.0₃0692'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.88'
$ws.Range('E23').Value = '  -5.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.11'
$ws.Range('E24').Value = '  -2.30%  '
$ws.Range('E25').Value = '  +4.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.22'
$ws.Range('E26').Value = '  -0.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.96'
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('E29').Value = '  -2.21%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.11'
$ws.Range('E31').Value = '  -3.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0464'
$ws.Range('E32').Value = '  -2.15%  '
$ws.Range('E33').Value = '  -0.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.14'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').Value = '1.391.67'
$ws.Range('E35').Value = '  -1.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.54'
$ws.Range('E36').Value = '  -2.93%  '
$ws.Range('E37').Value = '  -2.84%  '
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('E39').Value = '  +5.39%  '
$ws.Range('E40').Value = '  -0.93%  '
$ws.Range('E41').Value = '  -1.57%  '
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.793'
$ws.Range('E43').Value = '  -1.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.60'
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.85'
$ws.Range('E45').Value = '  +2.16%  '
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '63.00'
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('D48').Value = '1.709.68'
$ws.Range('E48').Value = '  -0.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.03'
$ws.Range('E49').Value = '  -1.29%  '
$ws.Range('D50').Value = '0.0₆0104'
$ws.Range('E50').Value = '  -1.79%  '
$ws.Range('E51').Value = '  -1.48%  '
